$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.257.98"
$ws.Range("E2").Value = "'  +1.58%  "
$ws.Range("D3").Value = "'1.646.24"
$ws.Range("E3").Value = "'  +0.43%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("D5").Value = "'217.59"
$ws.Range("E5").Value = "'  +0.86%  "
$ws.Range("E6").Value = "'  +0.27%  "
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E8").Value = "'  +0.15%  "
$ws.Range("E9").Value = "'  -0.05%  "
$ws.Range("D10").Value = "'20.04"
$ws.Range("E10").Value = "'  +1.28%  "
$ws.Range("D11").Value = "'0.0794"
$ws.Range("E11").Value = "'  +0.15%  "
$ws.Range("E12").Value = "'  +0.59%  "
$ws.Range("D13").Value = "'1.874.06"
$ws.Range("E13").Value = "'  +0.47%  "
$ws.Range("D14").Value = "'1.602.49"
$ws.Range("E14").Value = "'  -2.08%  "
$ws.Range("E15").Value = "'  -2.53%  "
$ws.Range("E16").Value = "'  -0.27%  "
$ws.Range("D17").Value = "'63.61"
$ws.Range("E17").Value = "'  +0.45%  "
$ws.Range("D18").Value = "'26.240.38"
$ws.Range("E18").Value = "'  +1.47%  "
$ws.Range("D20").Value = "'196.50"
$ws.Range("E20").Value = "'  +1.60%  "
$ws.Range("D21").Value = "'4.45"
$ws.Range("E21").Value = "'  -0.82%  "
$ws.Range("E23").Value = "'  -0.16%  "
$ws.Range("D24").Value = "'143.71"
$ws.Range("E24").Value = "'  +0.72%  "
$ws.Range("E25").Value = "'  -0.21%  "
$ws.Range("E26").Value = "'  -3.06%  "
$ws.Range("E27").Value = "'  +1.59%  "
$ws.Range("E28").Value = "'  -0.15%  "
$ws.Range("D29").Value = "'15.63"
$ws.Range("E29").Value = "'  +0.21%  "
$ws.Range("D30").Value = "'1.26"
$ws.Range("E30").Value = "'  +1.49%  "
$ws.Range("E31").Value = "'  +2.10%  "
$ws.Range("E32").Value = "'  -0.47%  "
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("E34").Value = "'  +1.28%  "
$ws.Range("E35").Value = "'  +0.84%  "
$ws.Range("E36").Value = "'  +0.35%  "
$ws.Range("D37").Value = "'1.137.23"
$ws.Range("E37").Value = "'  +0.28%  "
$ws.Range("E38").Value = "'  +1.59%  "
$ws.Range("E39").Value = "'  -1.79%  "
$ws.Range("E40").Value = "'  +0.53%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("D42").Value = "'5.69"
$ws.Range("E42").Value = "'  +2.21%  "
$ws.Range("D43").Value = "'100.18"
$ws.Range("E43").Value = "'  -0.34%  "
$ws.Range("D44").Value = "'0.798"
$ws.Range("E44").Value = "'  -1.35%  "
$ws.Range("D45").Value = "'1.782.93"
$ws.Range("E45").Value = "'  +0.47%  "
$ws.Range("D46").Value = "'56.42"
$ws.Range("E46").Value = "'  +1.66%  "
$ws.Range("E47").Value = "'  +3.49%  "
$ws.Range("E48").Value = "'  +2.47%  "
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'7.74"
$ws.Range("E49").Value = "'  +2.98%  "
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.417"
$ws.Range("E50").Value = "'  -0.09%  "
$ws.Range("D51").Value = "'0.0973"
$ws.Range("E51").Value = "'  +1.02%  "
